# Insert a new weekly price-observation row for "Betarraga" (Hortaliza,
# Macroferia Regional de Talca) at row 294, pushing the existing rows
# 294:419 down to 295:420.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(294).Insert()

$ws.Cells.Item(294, 1).Value = 5
$ws.Cells.Item(294, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(294, 3).Value = "Maule"
$ws.Cells.Item(294, 4).Value = 44875
$ws.Cells.Item(294, 5).Value = 7
$ws.Cells.Item(294, 6).Value = 100114014
$ws.Cells.Item(294, 7).Value = "Betarraga"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 5000
$ws.Cells.Item(294, 11).Value = 800
$ws.Cells.Item(294, 12).Value = 800
$ws.Cells.Item(294, 13).Value = 800
$ws.Cells.Item(294, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(294, 15).Value = "Región del Maule"
$ws.Cells.Item(294, 16).Value = 160
$ws.Cells.Item(294, 17).Value = 5
$ws.Cells.Item(294, 18).Value = "Hortaliza"
